$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = "x"
$ws.Range("B8").Value = 6.9
$ws.Range("C8").Value = "d"

$ws.Range("A9").Value = "y"
$ws.Range("B9").Value = 10.5
$ws.Range("C9").Value = "d"

$ws.Range("C9").Select()
